$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 168.76471
$ws.Range("I39").Value = 47.142857
$ws.Range("J39").Value = 736.3333
$ws.Range("K39").Value = 141.428571
$ws.Range("L39").Value = 2208.9999
$ws.Range("M39").Value = 154.571429
$ws.Range("N39").Value = -2800.9999
$ws.Range("H64").Value = 62502996
$ws.Range("I64").Value = 200001980
$ws.Range("J64").Value = 3460
$ws.Range("K64").Value = 200001980
$ws.Range("L64").Value = 3460
$ws.Range("M64").Value = -200001732
$ws.Range("N64").Value = -3956
$ws.Range("H67").Value = 62502996
$ws.Range("I67").Value = 200001980
$ws.Range("J67").Value = 3460
$ws.Range("K67").Value = 200001980
$ws.Range("L67").Value = 3460
$ws.Range("M67").Value = -200001122
$ws.Range("N67").Value = -5176
$ws.Range("H74").Value = 3102.658
$ws.Range("I74").Value = 2422.7856
$ws.Range("J74").Value = 3499.25
$ws.Range("K74").Value = 2422.7856
$ws.Range("L74").Value = 3499.25
$ws.Range("M74").Value = -1486.7856
$ws.Range("N74").Value = -5371.25
$ws.Range("H77").Value = 3102.658
$ws.Range("I77").Value = 2422.7856
$ws.Range("J77").Value = 3499.25
$ws.Range("K77").Value = 12113.928
$ws.Range("L77").Value = 17496.25
$ws.Range("M77").Value = -7433.928
$ws.Range("N77").Value = -26856.25
$ws.Range("H127").Value = 1215.3182
$ws.Range("I127").Value = 1168
$ws.Range("J127").Value = 1229.2354
$ws.Range("K127").Value = 3504
$ws.Range("L127").Value = 3687.7062
$ws.Range("M127").Value = 1456
$ws.Range("N127").Value = -13607.7062
$ws.Range("H137").Value = 3196.2144
$ws.Range("I137").Value = 1086.4242
$ws.Range("J137").Value = 6223.304
$ws.Range("K137").Value = 3259.2726
$ws.Range("L137").Value = 18669.912
$ws.Range("M137").Value = -709.2725999999998
$ws.Range("N137").Value = -23769.912
$ws.Range("H138").Value = 2447.386
$ws.Range("I138").Value = 1735.0938
$ws.Range("J138").Value = 3359.12
$ws.Range("K138").Value = 5205.2814
$ws.Range("L138").Value = 10077.36
$ws.Range("M138").Value = -65.28139999999985
$ws.Range("N138").Value = -20357.36
$ws.Range("H141").Value = 6295.5884
$ws.Range("I141").Value = 3385.484
$ws.Range("J141").Value = 36366.668
$ws.Range("K141").Value = 10156.452
$ws.Range("L141").Value = 109100.004
$ws.Range("M141").Value = -4976.451999999999
$ws.Range("N141").Value = -119460.004

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 218846.72
$ws.Range("I61").Value = 1509.2778
$ws.Range("J61").Value = 1001261.5
$ws.Range("K61").Value = 1509.2778
$ws.Range("L61").Value = 1001261.5
$ws.Range("M61").Value = -1297.2778
$ws.Range("N61").Value = -1001685.5
$ws.Range("H74").Value = 5294.074
$ws.Range("I74").Value = 1325.2142
$ws.Range("J74").Value = 9568.23
$ws.Range("K74").Value = 1325.2142
$ws.Range("L74").Value = 9568.23
$ws.Range("M74").Value = -451.2141999999999
$ws.Range("N74").Value = -11316.23
$ws.Range("H77").Value = 5294.074
$ws.Range("I77").Value = 1325.2142
$ws.Range("J77").Value = 9568.23
$ws.Range("K77").Value = 6626.071
$ws.Range("L77").Value = 47841.14999999999
$ws.Range("M77").Value = -2258.071
$ws.Range("N77").Value = -56577.14999999999
$ws.Range("H97").Value = 46707.637
$ws.Range("I97").Value = 91787.91
$ws.Range("J97").Value = 1627.3636
$ws.Range("K97").Value = 91787.91
$ws.Range("L97").Value = 1627.3636
$ws.Range("M97").Value = -91291.91
$ws.Range("N97").Value = -2619.3636
$ws.Range("H132").Value = 5131.094
$ws.Range("I132").Value = 3426.9048
$ws.Range("J132").Value = 11638
$ws.Range("K132").Value = 10280.7144
$ws.Range("L132").Value = 34914
$ws.Range("M132").Value = -7750.714399999999
$ws.Range("N132").Value = -39974
$ws.Range("H136").Value = 218846.72
$ws.Range("I136").Value = 1509.2778
$ws.Range("J136").Value = 1001261.5
$ws.Range("K136").Value = 4527.8334
$ws.Range("L136").Value = 3003784.5
$ws.Range("M136").Value = -1977.8334
$ws.Range("N136").Value = -3008884.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2696.6667
$ws.Range("I134").Value = 2000
$ws.Range("J134").Value = 2783.75
$ws.Range("K134").Value = 6000
$ws.Range("L134").Value = 8351.25
$ws.Range("M134").Value = -3465
$ws.Range("N134").Value = -13421.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5367.8887
$ws.Range("I16").Value = 2062.2
$ws.Range("J16").Value = 9500
$ws.Range("K16").Value = 2062.2
$ws.Range("L16").Value = 9500
$ws.Range("M16").Value = -1775.2
$ws.Range("N16").Value = -10074
$ws.Range("H31").Value = 3270.596
$ws.Range("I31").Value = 1141.159
$ws.Range("J31").Value = 4974.1455
$ws.Range("K31").Value = 1141.159
$ws.Range("L31").Value = 4974.1455
$ws.Range("M31").Value = -846.1590000000001
$ws.Range("N31").Value = -5564.1455
$ws.Range("H34").Value = 3270.596
$ws.Range("I34").Value = 1141.159
$ws.Range("J34").Value = 4974.1455
$ws.Range("K34").Value = 1141.159
$ws.Range("L34").Value = 4974.1455
$ws.Range("M34").Value = -939.1590000000001
$ws.Range("N34").Value = -5378.1455
$ws.Range("H86").Value = 3864.625
$ws.Range("I86").Value = 3100.2666
$ws.Range("J86").Value = 5138.5557
$ws.Range("K86").Value = 3100.2666
$ws.Range("L86").Value = 5138.5557
$ws.Range("M86").Value = -1977.2666
$ws.Range("N86").Value = -7384.5557
$ws.Range("H89").Value = 3864.625
$ws.Range("I89").Value = 3100.2666
$ws.Range("J89").Value = 5138.5557
$ws.Range("K89").Value = 15501.333
$ws.Range("L89").Value = 25692.7785
$ws.Range("M89").Value = -9885.332999999999
$ws.Range("N89").Value = -36924.7785
$ws.Range("H94").Value = 1690
$ws.Range("I94").Value = 1226.6666
$ws.Range("J94").Value = 2153.3333
$ws.Range("K94").Value = 1226.6666
$ws.Range("L94").Value = 2153.3333
$ws.Range("M94").Value = -775.6666
$ws.Range("N94").Value = -3055.3333
$ws.Range("H113").Value = 5367.8887
$ws.Range("I113").Value = 2062.2
$ws.Range("J113").Value = 9500
$ws.Range("K113").Value = 2062.2
$ws.Range("L113").Value = 9500
$ws.Range("M113").Value = 107.8000000000002
$ws.Range("N113").Value = -13840
$ws.Range("H134").Value = 5332.6
$ws.Range("I134").Value = 8722.4
$ws.Range("J134").Value = 1942.8
$ws.Range("K134").Value = 26167.2
$ws.Range("L134").Value = 5828.4
$ws.Range("M134").Value = -23632.2
$ws.Range("N134").Value = -10898.4
$ws.Range("H141").Value = 142911650
$ws.Range("I141").Value = 25000
$ws.Range("J141").Value = 166726080
$ws.Range("K141").Value = 25000
$ws.Range("L141").Value = 166726080
$ws.Range("M141").Value = -19820
$ws.Range("N141").Value = -166736440

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 685.725
$ws.Range("I122").Value = 283.08334
$ws.Range("J122").Value = 1289.6875
$ws.Range("K122").Value = 2547.75006
$ws.Range("L122").Value = 11607.1875
$ws.Range("M122").Value = -97.7500600000003
$ws.Range("N122").Value = -16507.1875
$ws.Range("H131").Value = 885.2449
$ws.Range("I131").Value = 327.33334
$ws.Range("J131").Value = 1010.775
$ws.Range("K131").Value = 982.0000200000001
$ws.Range("L131").Value = 3032.325
$ws.Range("M131").Value = 4057.99998
$ws.Range("N131").Value = -13112.325

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 37346.547
$ws.Range("I11").Value = 16759.2
$ws.Range("J11").Value = 54502.668
$ws.Range("K11").Value = 16759.2
$ws.Range("L11").Value = 54502.668
$ws.Range("M11").Value = -16620.2
$ws.Range("N11").Value = -54780.668
$ws.Range("H12").Value = 2176.6428
$ws.Range("I12").Value = 2176.6428
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 2176.6428
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -2036.6428
$ws.Range("H13").Value = 236.2
$ws.Range("I13").Value = 100
$ws.Range("J13").Value = 270.25
$ws.Range("K13").Value = 100
$ws.Range("L13").Value = 270.25
$ws.Range("M13").Value = 39
$ws.Range("N13").Value = -548.25
$ws.Range("H132").Value = 19791.402
$ws.Range("I132").Value = 30166.795
$ws.Range("J132").Value = 2198.348
$ws.Range("K132").Value = 90500.38499999999
$ws.Range("L132").Value = 6595.044
$ws.Range("M132").Value = -87970.38499999999
$ws.Range("N132").Value = -11655.044

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 1006
$ws.Range("I13").Value = 1006
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 1006
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -866
$ws.Range("H20").Value = 10000
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 10000
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 10000
$ws.Range("N20").Value = -10452
$ws.Range("H99").Value = 30252.666
$ws.Range("I99").Value = 30879
$ws.Range("J99").Value = 29000
$ws.Range("K99").Value = 30879
$ws.Range("L99").Value = 29000
$ws.Range("M99").Value = -27884
$ws.Range("N99").Value = -34990

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1857.2727
$ws.Range("I126").Value = 2191.6667
$ws.Range("J126").Value = 1456
$ws.Range("K126").Value = 6575.000100000001
$ws.Range("L126").Value = 4368
$ws.Range("M126").Value = -4105.000100000001
$ws.Range("N126").Value = -9308
$ws.Range("H132").Value = 2999.8208
$ws.Range("I132").Value = 3251.827
$ws.Range("J132").Value = 2126.2
$ws.Range("K132").Value = 9755.481
$ws.Range("L132").Value = 6378.599999999999
$ws.Range("M132").Value = -7225.481
$ws.Range("N132").Value = -11438.6
